$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (make_model) to hold the "year" field
$ws.Range("D1").EntireColumn.Insert()

# Header for new column
$ws.Range("D1").Value = "year"

# Year values for each aircraft row (rows 2-11)
$years = @(1982, 1968, 1977, 1984, 1967, 1959, 1963, 1968, 1970, 1960)
for ($i = 0; $i -lt $years.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $years[$i]
}

# Update selection to match the post-edit active cell
$ws.Range("D12").Select()
